$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Plain numeric value updates (same cell type before/after) ---
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 122.222222222222
$ws.Range("M15").Value = 53.846153846153
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 119
$ws.Range("K16").Value = -27.731092436974
$ws.Range("L16").Value = -37.681159420289
$ws.Range("M16").Value = -43.790849673202
$ws.Range("N16").Value = -87.955182072829
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 166
$ws.Range("K17").Value = -7.228915662650
$ws.Range("L17").Value = -8.875739644970
$ws.Range("M17").Value = 52.475247524752
$ws.Range("N17").Value = -44.404332129963
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 250
$ws.Range("I18").Value = 50
$ws.Range("K18").Value = 6.382978723404
$ws.Range("L18").Value = -43.181818181818
$ws.Range("M18").Value = -31.506849315068
$ws.Range("N18").Value = -91.039426523297
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -54.545454545454
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 347
$ws.Range("J19").Value = 329
$ws.Range("K19").Value = 5.471124620060
$ws.Range("L19").Value = -3.878116343490
$ws.Range("M19").Value = 39.357429718875
$ws.Range("N19").Value = -44.301765650080
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 2.380952380952
$ws.Range("L20").Value = -57
$ws.Range("N20").Value = -87.240356083086
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 1.639344262295
$ws.Range("I21").Value = 702
$ws.Range("J21").Value = 713
$ws.Range("K21").Value = -1.542776998597
$ws.Range("L21").Value = -18.843930635838
$ws.Range("M21").Value = 13.776337115072
$ws.Range("N21").Value = -72.427336999214
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("I22").Value = 24
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = -7.692307692307
$ws.Range("L22").Value = -29.411764705882
$ws.Range("M22").Value = -22.580645161290
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 9.090909090909
$ws.Range("I23").Value = 102
$ws.Range("J23").Value = 122
$ws.Range("K23").Value = -16.393442622950
$ws.Range("L23").Value = -27.659574468085
$ws.Range("M23").Value = 39.726027397260
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 37.5
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 18.918918918918
$ws.Range("I24").Value = 411
$ws.Range("J24").Value = 380
$ws.Range("K24").Value = 8.157894736842
$ws.Range("L24").Value = -11.991434689507
$ws.Range("M24").Value = -27.640845070422
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 75
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 7.142857142857
$ws.Range("L25").Value = -48.979591836734
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -54.166666666666
$ws.Range("I26").Value = 230
$ws.Range("J26").Value = 228
$ws.Range("K26").Value = 0.877192982456
$ws.Range("L26").Value = 3.139013452914
$ws.Range("M26").Value = -14.498141263940
$ws.Range("D27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 53.333333333333
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = -17.948717948717
$ws.Range("L28").Value = -27.272727272727
$ws.Range("N29").Value = -84.375
$ws.Range("N30").Value = -81.481481481481
$ws.Range("L31").Value = -38.461538461538

# --- Type-change cells (number <-> text placeholder) ---
# Style 13 = right-aligned text cell (used for "0"/"***.*" placeholders).
# Style 14 = right-aligned integer cell (#,##0).
# $ws.Range("C15") is a stable style-13 anchor; $ws.Range("I15") is a stable style-14 anchor.
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("C17").NumberFormat = "General"
$ws.Range("C17").Value = 4
$ws.Range("I15").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("C18").NumberFormat = "General"
$ws.Range("C18").Value = 2
$ws.Range("I15").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("C25").NumberFormat = "General"
$ws.Range("C25").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F27").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F31").PasteSpecial(-4122)

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = $false